# Slide 5: reposition "Picture 8" and "Picture 10".
# Target offsets (EMU): Picture 8 -> (943256, 3469248); Picture 10 -> (709370, 2630696)
# Shape.Left/.Top are expressed in points (1 pt = 12700 EMU) and are stored as
# single-precision floats internally, so the literals below were chosen (via a
# small nearby-float search) to round-trip to the exact target EMU values
# rather than using the naive EMU/12700 quotient (which can truncate down by
# 1 EMU once cast to a 32-bit float).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$pic8 = $s.Shapes.Item("Picture 8")
$pic8.Left = 74.27213287353516
$pic8.Top = 273.1691589355469

$pic10 = $s.Shapes.Item("Picture 10")
$pic10.Left = 55.85590744018555
$pic10.Top = 207.14141845703125
